$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 14 -----
$ws.Range("A14").Value = 11
$ws.Range("C14").Value = "contains duplicates 3"
$ws.Range("B14").Value = "todo"
$ws.Range("D14").Value = "Medium"
$ws.Range("E14").Value = "bst"
$ws.Range("F14").Value = "Medium"
$ws.Range("H14").Value = 60

# ----- Row 15 -----
$ws.Range("A15").Value = 12
$ws.Range("C15").Value = "Merge k sorted list"
$ws.Range("B15").Value = "todo 23"
$ws.Range("D15").Value = "Hard"

# ----- Row 16 -----
$ws.Range("A16").Value = 12
$ws.Range("B16").Value = "todo 41"
$ws.Range("C16").Value = "first missing positive"
$ws.Range("D16").Value = "Hard"

# Update selection to match the resulting active cell
$ws.Range("F16").Select()
